{"js": "// Apply strikethrough formatting to several character/paragraph ranges, as\n// described by the target diff:\n//  1. \"She has two children, Stephanie and Shawn, who are currently\n//     pursuing their studies.\" (tail of the Teng Yan description paragraph)\n//  2. The whole \"Stephanie Lee (affectionately called Ah Cheh...)\" paragraph\n//  3. The whole \"Teng Yan's daughter, who has a young sibling Shawn...\" paragraph\n//  4. The whole \"Shawn Lee (affectionately known as Ah Zai...)\" paragraph\n//  5. The whole \"Teng Yan's son, who is currently studying...\" paragraph\n//  6. Just the word \"members\" inside \"Friction between the family members...\"\n//  7. The whole \"Day care centre\" paragraph\n//  8. The whole \"Hoping to engage Ah Ma in more activities...\" paragraph\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfunction findParagraph(startsWith) {\n  for (const p of paragraphs.items) {\n    if (p.text.indexOf(startsWith) === 0) {\n      return p;\n    }\n  }\n  throw new Error(\"Paragraph not found: \" + startsWith);\n}\n\n// 1. Split run: strike only the second sentence about her two children.\nconst tengYanDescription = findParagraph(\"Ah Ma\\u2019s elder daughter and her primary caregiver.\");\nconst childrenSentence = tengYanDescription.search(\n  \"She has two children, Stephanie and Shawn, who are currently pursuing their studies.\",\n  { matchCase: true }\n);\nawait context.sync();\nchildrenSentence.items[0].font.strikeThrough = true;\n\n// 2-5. Strike entire paragraphs (run text + paragraph mark).\nfindParagraph(\"Stephanie Lee (affectionately called Ah Cheh\").font.strikeThrough = true;\nfindParagraph(\"Teng Yan\\u2019s daughter, who has a young sibling Shawn\").font.strikeThrough = true;\nfindParagraph(\"Shawn Lee (affectionately known as Ah Zai\").font.strikeThrough = true;\nfindParagraph(\"Teng Yan\\u2019s son, who is currently studying in secondary school\").font.strikeThrough = true;\n\n// 6. Strike just the word \"members\" within the friction paragraph.\nconst frictionParagraph = findParagraph(\"Friction between the family members\");\nconst membersWord = frictionParagraph.search(\"members\", { matchCase: true });\nawait context.sync();\nmembersWord.items[0].font.strikeThrough = true;\n\n// 7-8. Strike entire paragraphs (run text + paragraph mark).\nfindParagraph(\"Day care centre\").font.strikeThrough = true;\nfindParagraph(\"Hoping to engage Ah Ma in more activities\").font.strikeThrough = true;\n\nawait context.sync();\n", "ps1": "# Apply strikethrough formatting to several character/paragraph ranges, as\n# described by the target diff:\n#  1. \"She has two children, Stephanie and Shawn, who are currently\n#     pursuing their studies.\" (tail of the Teng Yan description paragraph)\n#  2. The whole \"Stephanie Lee (affectionately called Ah Cheh...)\" paragraph\n#  3. The whole \"Teng Yan's daughter, who has a young sibling Shawn...\" paragraph\n#  4. The whole \"Shawn Lee (affectionately known as Ah Zai...)\" paragraph\n#  5. The whole \"Teng Yan's son, who is currently studying...\" paragraph\n#  6. Just the word \"members\" inside \"Friction between the family members...\"\n#  7. The whole \"Day care centre\" paragraph\n#  8. The whole \"Hoping to engage Ah Ma in more activities...\" paragraph\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphStartingWith($prefix) {\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text.StartsWith($prefix)) {\n            return $p\n        }\n    }\n    throw \"Paragraph not found: $prefix\"\n}\n\n# 1. Split run: strike only the second sentence about her two children.\n$tengYanDescription = Get-ParagraphStartingWith \"Ah Ma\u2019s elder daughter and her primary caregiver.\"\n$childrenRange = $tengYanDescription.Range\n$childrenRange.Find.Execute(\"She has two children, Stephanie and Shawn, who are currently pursuing their studies.\") | Out-Null\n$childrenRange.Font.StrikeThrough = 1\n\n# 2-5. Strike entire paragraphs (run text + paragraph mark).\n(Get-ParagraphStartingWith \"Stephanie Lee (affectionately called Ah Cheh\").Range.Font.StrikeThrough = 1\n(Get-ParagraphStartingWith \"Teng Yan\u2019s daughter, who has a young sibling Shawn\").Range.Font.StrikeThrough = 1\n(Get-ParagraphStartingWith \"Shawn Lee (affectionately known as Ah Zai\").Range.Font.StrikeThrough = 1\n(Get-ParagraphStartingWith \"Teng Yan\u2019s son, who is currently studying in secondary school\").Range.Font.StrikeThrough = 1\n\n# 6. Strike just the word \"members\" within the friction paragraph.\n$frictionParagraph = Get-ParagraphStartingWith \"Friction between the family members\"\n$membersRange = $frictionParagraph.Range\n$membersRange.Find.Execute(\"members\") | Out-Null\n$membersRange.Font.StrikeThrough = 1\n\n# 7-8. Strike entire paragraphs (run text + paragraph mark).\n(Get-ParagraphStartingWith \"Day care centre\").Range.Font.StrikeThrough = 1\n(Get-ParagraphStartingWith \"Hoping to engage Ah Ma in more activities\").Range.Font.StrikeThrough = 1\n"}
